# Resolve indentation problem on build_model.py add Model_7 for ID = 1745
#
# Adds a new table row (row 14) for Model_7 / Speed ID 1745, matching
# the existing table's data + formatting conventions, and nudges a
# couple of cosmetic workbook settings (theme accent colors, selection)
# that travelled along with the author's save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New data row (row 14) -------------------------------------------
# MODEL ID
$ws.Range("A14").Value = 7
# SPEED ID
$ws.Range("B14").Value = 1745
# INPUTS
$ws.Range("C14").WrapText = $true
$ws.Range("C14").Value = "scaled speed`nweekday o.h.`ndaypart o.h."
# NN LAYERS
$ws.Range("D14").WrapText = $true
$ws.Range("D14").Value = "lstm(50)+do(.3)`nlstm/50)+do(.3)`nlstm/33)"
# EPOCH
$ws.Range("E14").Value = 100
# TIME STAMP
$ws.Range("F14").WrapText = $true
$ws.Range("F14").Value = "1h back`n1h forward"
# ESTIMATION DISTANCE
$ws.Range("G14").WrapText = $true
$ws.Range("G14").Value = "1 week+`n2 week+`n3 week"
# TRAIN DATA
$ws.Range("H14").WrapText = $true
$ws.Range("H14").Value = "Feb March April May"
# TEST DATA
$ws.Range("I14").WrapText = $true
$ws.Range("I14").Value = "First 7 days of June"
# TRAIN ERROR
$ws.Range("J14").WrapText = $true
$ws.Range("J14").Value = 20.11
# TEST ERROR
$ws.Range("K14").WrapText = $true
$ws.Range("K14").Value = 37.41
# RUSH ER.
$ws.Range("L14").Value = 59.076
# COMMENT
$ws.Range("M14").WrapText = $true
$ws.Range("M14").Value = "Adding 3 week before the estimation hour decreased the success instead of increasing it. Why test loss consistently rises="

$ws.Rows.Item(14).RowHeight = 33

# ---- Theme tweak that rode along with this save -----------------------
$theme = $wb.Theme
$colorScheme = $theme.ThemeColorScheme
$colorScheme.Colors(1).RGB = 65535   # dk1: 000000 -> FFFF00
$colorScheme.Colors(2).RGB = 0       # lt1: FFFFFF -> 000000

# Legacy indexed color palette touched by the save as well.
$wb.ResetColors()

# ---- Selection / scroll position on exit -------------------------------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("M15").Select()

Write-Output "Model_7 (Speed ID 1745) row added"
